$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the closing-row dates (B9/C9) that were previously blank, matching
# the "d-mmm" date format already used by the rows above (B5:C8), which keep
# a medium border on top and bottom (same as the header row style).
$ws.Range("B9").NumberFormat = "d-mmm"
$ws.Range("C9").NumberFormat = "d-mmm"
$ws.Range("B9").Value = 45531
$ws.Range("C9").Value = 45531

# Update the active selection to A4:C4
$ws.Range("A4:C4").Select() | Out-Null
